$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original wide-format data (row1 headers B1:G1, row2 values B2:G2)
$headers = @()
$values = @()
for ($col = 2; $col -le 7; $col++) {
    $headers += $ws.Cells.Item(1, $col).Value2
    $values += $ws.Cells.Item(2, $col).Value2
}

# Clear existing contents so we can rebuild the sheet in long format
$ws.Cells.Clear()

# New header row: A1 blank, B1 = Condition, C1 = ddCT
$ws.Range("B1").Value = "Condition"
$ws.Range("C1").Value = "ddCT"

$styledRange = $ws.Range("B1:C1")
$styledRange.Font.Bold = $true
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160
$styledRange.Borders.LineStyle = 1

for ($i = 0; $i -lt $headers.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $headers[$i]
    $ws.Cells.Item($r, 3).Value = $values[$i]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
}
